# 11:26 time 01.11.2024 date
# Append 8 new applicant rows (88-95) to the qabul sheet. All columns in
# this sheet are stored as plain text in the source file (even the
# contract numbers, phone numbers and dd-mm-yyyy dates), so each row is
# written with NumberFormat forced to "@" first to stop Excel's COM
# layer from auto-converting numeric-looking / date-looking strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    [PSCustomObject]@{ Row=88; A="Aliyeva Mahliyo Murodjon qizi"; B="Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"; C="AD8955000"; D="684"; E="Toshkent viloyati"; F="Piskent tumani"; G="998953500128"; H="31-10-2024" }
    [PSCustomObject]@{ Row=89; A="Saxtabova Umidaxon Bekxo'ja qizi"; B="Defektologiya (logopediya) 576 soatlik"; C="AD6680957"; D="685"; E="Toshkent viloyati"; F="Yangiyoʻl tumani"; G="998930486446"; H="31-10-2024" }
    [PSCustomObject]@{ Row=90; A="Qipchakova Etiborxon Yuldashevna"; B="Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"; C="AD8410114"; D="686"; E="Fargona viloyati"; F="Buvayda tumani"; G="998911440526"; H="31-10-2024" }
    [PSCustomObject]@{ Row=91; A="Djaborova Marhaboxon Zokirjonovna"; B="Defektologiya (logopediya) 576 soatlik"; C="AD1322435"; D="687"; E="Toshkent viloyati"; F="Ohangaron tumani"; G="+998931893234"; H="31-10-2024" }
    [PSCustomObject]@{ Row=92; A="Axtamova Marjona Utkir qizi"; B="Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"; C="AC0716066"; D="688"; E="Navoiy viloyati"; F="Navbahor tumani"; G="998505064800"; H="31-10-2024" }
    [PSCustomObject]@{ Row=93; A="Zokirova Soliha Abduraxim qizi"; B="Defektologiya (logopediya) 576 soatlik"; C="AD0626175"; D="689"; E="Toshkent shahri"; F="Shayxontohur tumani"; G="998974335353"; H="31-10-2024" }
    [PSCustomObject]@{ Row=94; A="Yo'ldosheva Mohlaroyim To'ymurod qizi"; B="Defektologiya (logopediya) 576 soatlik"; C="AD0130593"; D="690"; E="Navoiy viloyati"; F="Qiziltepa tumani"; G="998999264272"; H="01-11-2024" }
    [PSCustomObject]@{ Row=95; A="Rahmonova Muxlisa Gofurgon qizi"; B="Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik"; C="AB2555334"; D="691"; E="Samarqand viloyati"; F="Narpay tumani"; G="998978927477"; H="01-11-2024" }
)

foreach ($r in $newRows) {
    $rowRange = $ws.Range("A" + $r.Row + ":H" + $r.Row)

    # Force text storage before assigning values.
    $rowRange.NumberFormat = "@"

    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("H" + $r.Row).Value = $r.H

    # Clear the temporary number-format again so the new cells keep the
    # same "no explicit style" look as the rest of the data rows.
    $rowRange.ClearFormats()
}
